$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 468, shifting existing rows 468:500 down to 469:501
$ws.Rows("468:468").Insert()

# Populate the newly inserted row 468 with the new record
$ws.Range("A468").Value = 5
$ws.Range("B468").Value = "Macroferia Regional de Talca"
$ws.Range("C468").Value = "Maule"
$ws.Range("D468").Value = 45106
$ws.Range("E468").Value = 7
$ws.Range("F468").Value = 100112003
$ws.Range("G468").Value = "Ajo"
$ws.Range("H468").Value = "Chino"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 300
$ws.Range("K468").Value = 18000
$ws.Range("L468").Value = 18000
$ws.Range("M468").Value = 18000
$ws.Range("N468").Value = "`$/malla 10 kilos"
$ws.Range("O468").Value = "China"
$ws.Range("P468").Value = 1800
$ws.Range("Q468").Value = 10
$ws.Range("R468").Value = "Hortaliza"
